# Regenerate save_data: column G ("K" = strikes/putts, formerly "Strike#")
# is recomputed from the freshly calculated std/mean "s_vals". Apply the
# newly computed K values to the existing rows on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 0
    3  = 0
    4  = 0
    5  = 1
    7  = 2
    8  = 2
    9  = 1
    10 = 0
    11 = 1
    12 = 0
    13 = 1
    14 = 2
    15 = 0
    16 = 0
    17 = 1
    18 = 1
    19 = 1
    20 = 2
    21 = 1
    22 = 0
    24 = 0
    25 = 0
    26 = 2
    27 = 1
    28 = 2
    29 = 1
    30 = 1
    31 = 1
    32 = 0
    33 = 1
    35 = 1
    36 = 2
    37 = 1
    38 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
